$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to store as text (preserve formatting like trailing
# zeros / multi-dot "thousands" strings) without leaving the cells style
# changed from the workbook default.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '30.454.43'
$ws.Range('E2').Value = '  +0.46%  '
Set-TextValue $ws.Range('D3') '2.106.00'
$ws.Range('E3').Value = '  +4.76%  '
Set-TextValue $ws.Range('D4') '1.002'
$ws.Range('E4').Value = '  -0.01%  '
Set-TextValue $ws.Range('D5') '329.33'
$ws.Range('E5').Value = '  +1.48%  '
Set-TextValue $ws.Range('D6') '1.001'
$ws.Range('E6').Value = '  -0.02%  '
Set-TextValue $ws.Range('D7') '0.5256'
$ws.Range('E7').Value = '  +2.53%  '
Set-TextValue $ws.Range('D8') '0.4361'
$ws.Range('E8').Value = '  +2.06%  '
Set-TextValue $ws.Range('D9') '0.08860'
$ws.Range('E9').Value = '  +1.95%  '
Set-TextValue $ws.Range('D10') '47.19'
$ws.Range('E10').Value = '  +9.47%  '
Set-TextValue $ws.Range('D11') '1.163'
$ws.Range('E11').Value = '  +2.46%  '
Set-TextValue $ws.Range('D12') '24.51'
$ws.Range('E12').Value = '  -0.83%  '
Set-TextValue $ws.Range('D13') '2.104.13'
$ws.Range('E13').Value = '  +4.77%  '
Set-TextValue $ws.Range('D14') '6.734'
$ws.Range('E14').Value = '  +2.75%  '
Set-TextValue $ws.Range('D15') '7.769'
$ws.Range('E15').Value = '  +4.15%  '
Set-TextValue $ws.Range('D16') '96.50'
$ws.Range('E16').Value = '  +2.32%  '
Set-TextValue $ws.Range('D17') '1.002'
$ws.Range('E17').Value = '  -0.02%  '
Set-TextValue $ws.Range('D18') '0.00001128'
$ws.Range('E18').Value = '  +1.18%  '
Set-TextValue $ws.Range('D19') '0.06639'
$ws.Range('E19').Value = '  +1.81%  '
Set-TextValue $ws.Range('D20') '19.01'
$ws.Range('E20').Value = '  +0.67%  '
Set-TextValue $ws.Range('D21') '1.001'
$ws.Range('E21').Value = '  -0.05%  '
Set-TextValue $ws.Range('D22') '6.337'
$ws.Range('E22').Value = '  +2.33%  '
Set-TextValue $ws.Range('D23') '30.517.20'
$ws.Range('E23').Value = '  +0.45%  '
Set-TextValue $ws.Range('D24') '12.33'
$ws.Range('E24').Value = '  +4.33%  '
Set-TextValue $ws.Range('D25') '2.335'
$ws.Range('E25').Value = '  +3.98%  '
Set-TextValue $ws.Range('D26') '2.350.53'
$ws.Range('E26').Value = '  +4.70%  '
Set-TextValue $ws.Range('D27') '22.41'
$ws.Range('E27').Value = '  -0.21%  '
Set-TextValue $ws.Range('D28') '2.586'
$ws.Range('E28').Value = '  +6.70%  '
Set-TextValue $ws.Range('D29') '161.92'
$ws.Range('E29').Value = '  -0.23%  '
Set-TextValue $ws.Range('D30') '132.70'
$ws.Range('E30').Value = '  +1.27%  '
$ws.Range('E31').Value = '  +6.80%  '
$ws.Range('E32').Value = '  +2.40%  '
Set-TextValue $ws.Range('D33') '1.680'
$ws.Range('E33').Value = '  +22.69%  '
Set-TextValue $ws.Range('D34') '6.206'
$ws.Range('E34').Value = '  +2.36%  '
Set-TextValue $ws.Range('D35') '3.920'
$ws.Range('E35').Value = '  +2.43%  '
Set-TextValue $ws.Range('D36') '9.980'
$ws.Range('E36').Value = '  +9.81%  '
$ws.Range('E37').Value = '  +2.56%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D38') '0.06704'
$ws.Range('E38').Value = '  +0.66%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D39') '5.484'
$ws.Range('E39').Value = '  +0.35%  '
Set-TextValue $ws.Range('D40') '12.68'
$ws.Range('E40').Value = '  +2.69%  '
Set-TextValue $ws.Range('D41') '0.2278'
$ws.Range('E41').Value = '  +4.02%  '
Set-TextValue $ws.Range('D42') '0.6818'
$ws.Range('E42').Value = '  +2.93%  '
Set-TextValue $ws.Range('D43') '1.260'
$ws.Range('E43').Value = '  +2.14%  '
Set-TextValue $ws.Range('D44') '1.000'
$ws.Range('E44').Value = '  -0.03%  '
Set-TextValue $ws.Range('D45') '14.02'
$ws.Range('E45').Value = '  +3.38%  '
Set-TextValue $ws.Range('D46') '0.6388'
$ws.Range('E46').Value = '  +3.64%  '
Set-TextValue $ws.Range('D47') '2.210'
$ws.Range('E47').Value = '  +1.47%  '
Set-TextValue $ws.Range('D48') '3.625'
$ws.Range('E48').Value = '  -0.88%  '
Set-TextValue $ws.Range('D49') '1.253'
$ws.Range('E49').Value = '  -0.67%  '
Set-TextValue $ws.Range('D50') '1.197'
$ws.Range('E50').Value = '  +8.20%  '
Set-TextValue $ws.Range('D51') '82.60'
$ws.Range('E51').Value = '  +2.57%  '
